# "fix bug compare two objects" -- add a PlaceHolder column (D) to Sheet1
# that holds a plain 0/1-style flag for each row, instead of relying on
# comparing two objects directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell for the extra column.
$ws.Range("D1").Value = "PlaceHolder"

# Seed the new column with a default (non-object) comparison flag per row.
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0

# Restore "portrait" page setup for the sheet (adds <pageSetup .../>).
$ws.PageSetup.Orientation = 1

# Leave the selection where the fix was made.
$ws.Range("D6").Select()
